$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "24.661.87"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.695.29"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.35%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "315.00"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3916"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4050"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "1.499"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.75%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "53.00"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.73%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08764"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.97%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "7.658"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +5.39%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "24.49"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.47%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.00001364"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.34%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "7.984"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.692.12"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "98.51"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.07118"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "19.83"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.362"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +4.43%  "

$ws.Range("E22").Value = "  +0.57%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "14.30"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "24.647.85"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.027"
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.355"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  -0.27%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "162.74"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.54%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.468"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +12.75%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "137.43"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.242"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.878.59"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.28%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.08922"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.80%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "7.535"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +5.62%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.050"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.99%  "

$ws.Range("E36").Value = "  +4.21%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.02941"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +7.82%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2737"
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "10.80"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -5.17%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "14.31"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.09130"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.7918"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.467"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "16.74"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +4.22%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.7230"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.577"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.36%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.221"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.27%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.330"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "139.18"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "91.23"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.14%  "
